$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''48.112.81'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").Value = '''2.502.00'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''320.28'
$ws.Range("E5").Value = '  -0.83%  '
$ws.Range("D6").Value = '''107.44'
$ws.Range("E6").Value = '  -1.36%  '
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").Value = '''0.542'
$ws.Range("E9").Value = '  -1.54%  '
$ws.Range("D10").Value = '''39.73'
$ws.Range("E10").Value = '  -1.04%  '
$ws.Range("E11").Value = '  +6.22%  '
$ws.Range("D12").Value = '''0.0811'
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("E13").Value = '  -0.01%  '
$ws.Range("D15").Value = '''2.893.05'
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").Value = '''2.503.76'
$ws.Range("E16").Value = '  +0.36%  '
$ws.Range("D17").Value = '''0.835'
$ws.Range("E17").Value = '  -1.75%  '
$ws.Range("D18").Value = '''47.979.82'
$ws.Range("E18").Value = '  +0.67%  '
$ws.Range("D19").Value = '''12.97'
$ws.Range("E19").Value = '  -1.39%  '
$ws.Range("D20").Value = '''6.68'
$ws.Range("E20").Value = '  +0.81%  '
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("D23").Value = '''276.31'
$ws.Range("E23").Value = '  +11.50%  '
$ws.Range("D24").Value = '''71.51'
$ws.Range("E24").Value = '  +1.10%  '
$ws.Range("D25").Value = '''2.54'
$ws.Range("E25").Value = '  -0.60%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("D28").Value = '''2.29'
$ws.Range("E28").Value = '  +4.46%  '
$ws.Range("D29").Value = '''9.72'
$ws.Range("E29").Value = '  -2.54%  '
$ws.Range("E30").Value = '  +1.92%  '
$ws.Range("D31").Value = '''35.12'
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("D32").Value = '''49.72'
$ws.Range("E32").Value = '  -0.41%  '
$ws.Range("D33").Value = '''19.47'
$ws.Range("E33").Value = '  -1.98%  '
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").Value = '''5.30'
$ws.Range("E35").Value = '  -0.89%  '
$ws.Range("E36").Value = '  -0.72%  '
$ws.Range("E37").Value = '  -0.84%  '
$ws.Range("D38").Value = '''4.63'
$ws.Range("E38").Value = '  -0.86%  '
$ws.Range("D39").Value = '''2.88'
$ws.Range("E39").Value = '  -2.84%  '
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("D41").Value = '''121.26'
$ws.Range("E41").Value = '  +1.67%  '
$ws.Range("E42").Value = '  -0.26%  '
$ws.Range("D43").Value = '''21.47'
$ws.Range("E43").Value = '  -3.67%  '
$ws.Range("E44").Value = '  +1.92%  '
$ws.Range("D45").Value = '''2.020.60'
$ws.Range("E45").Value = '  +0.93%  '
$ws.Range("D46").Value = '''3.14'
$ws.Range("E46").Value = '  +2.73%  '
$ws.Range("D47").Value = '''1.99'
$ws.Range("E47").Value = '  -1.65%  '
$ws.Range("E48").Value = '  +1.81%  '
$ws.Range("E49").Value = '  -0.37%  '
$ws.Range("D50").Value = '''5.17'
$ws.Range("E50").Value = '  +1.16%  '
$ws.Range("D51").Value = '''80.39'
$ws.Range("E51").Value = '  +3.45%  '
